# Actualización automática 2025-11-19 10:30:09
# A new asesor ("MUÑOZ CALDERON JUAN ADOLFO") is inserted alphabetically
# (between MOROCHO PLAZA SHIRLEY AURELIA and PALMA PICO OSCAR FILIDEL)
# into both the "VENTAS POR GRUPO" and "VENTA MENSUAL" sheets. This pushes
# every following row down by one and bumps the "X de 54" -> "X de 55"
# summary counters at the bottom of "VENTAS POR GRUPO".

$wb = $excel.ActiveWorkbook

# ---- Sheet 1: VENTAS POR GRUPO (columns A:R, data rows 2..55 after insert) ----
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Insert a new row before row 42 (currently "PALMA PICO OSCAR FILIDEL"),
# shifting rows 42..56 down to 43..57.
$ws1.Rows.Item(42).Insert()

# Fill in the new asesor row with zeroed sales figures.
$ws1.Cells.Item(42, 1).Value = "OFICINA-CATAECSA"
$ws1.Cells.Item(42, 2).Value = "MUÑOZ CALDERON JUAN ADOLFO"
for ($col = 3; $col -le 18; $col++) {
    $ws1.Cells.Item(42, $col).Value = 0
}

# Update the "X de 54" -> "X de 55" summary counters now on row 57.
for ($col = 3; $col -le 18; $col++) {
    $cell = $ws1.Cells.Item(57, $col)
    $cell.Value = $cell.Value2 -replace "de 54", "de 55"
}

# ---- Sheet 2: VENTA MENSUAL (columns A:G, data rows 2..55 after insert) ----
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

# Insert a new row before row 42 (currently "PALMA PICO OSCAR FILIDEL"),
# shifting rows 42..56 down to 43..57.
$ws2.Rows.Item(42).Insert()

# Fill in the new asesor row with zeroed sales figures.
$ws2.Cells.Item(42, 1).Value = "OFICINA-CATAECSA"
$ws2.Cells.Item(42, 2).Value = "MUÑOZ CALDERON JUAN ADOLFO"
for ($col = 3; $col -le 7; $col++) {
    $ws2.Cells.Item(42, $col).Value = 0
}
